$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 10 so that Animal 163 (currently only row 9)
# gets three rows total (9, 10, 11), matching the pattern used for Animal 121.
# This shifts the old rows 10-19 down to 12-21.
$ws.Rows("10:11").Insert()

# --- Update Animal 121 rows (6-8): Start Time / Window # values ---
$ws.Range("D6").Value = 0.45925925925925926
$ws.Range("E6").Value = 4

$ws.Range("D7").Value = 0.46249999999999997
$ws.Range("E7").Value = 3

$ws.Range("D8").Value = 0.46574074074074073
$ws.Range("E8").Value = 2

# --- Animal 163 row 9 (existing row, values updated) ---
$ws.Range("D9").Value = 0.45925925925925926
$ws.Range("E9").Value = 4

# --- Animal 163 new rows 10-11 (fill in data for the newly inserted rows) ---
$ws.Range("A10").Value = 163
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0.46249999999999997
$ws.Range("E10").Value = 3

$ws.Range("A11").Value = 163
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 0.46574074074074073
$ws.Range("E11").Value = 2

# --- Update the selection to match the new state ---
$ws.Range("E11").Select()
